$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Replace comma-separated category label patterns with slash-separated versions
$null = $ws.Cells.Replace("Lakossági összetétel, közösségek", "Lakossági összetétel / közösségek")
$null = $ws.Cells.Replace("etnikai, osztálybeli, stb.)", "etnikai / osztálybeli / más)")
$null = $ws.Cells.Replace("zaj, szmog, stb)", "zaj / szmog / más)")

# 2. Fix document filename typos
$null = $ws.Cells.Replace("PCSMZ02.docx", "PCMZ02.docx")
$null = $ws.Cells.Replace("PCSSA01.docx", "PCSA01.docx")
$null = $ws.Cells.Replace("PCSSA02.docx", "PCSA02.docx")
$null = $ws.Cells.Replace("PCSMZ01.docx", "PCMZ01.docx")
$null = $ws.Cells.Replace("DBSZA01.docx", "DBSA01.docx")

# 3. Swap rows 337 and 338 quotation/codes (fix misaligned code assignment)
$b337 = $ws.Cells.Item(337,2).Value()
$b338 = $ws.Cells.Item(338,2).Value()
$c337 = $ws.Cells.Item(337,3).Value()
$c338 = $ws.Cells.Item(338,3).Value()

$ws.Cells.Item(337,2).Value2 = $b338
$ws.Cells.Item(338,2).Value2 = $b337
$ws.Cells.Item(337,3).Value2 = $c338
$ws.Cells.Item(338,3).Value2 = $c337
